# Update "想去人数" (interested count) values per commit 456a3b4 gh-pages data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 5
$ws.Range("F8").Value = 42
$ws.Range("F9").Value = 725
$ws.Range("F10").Value = 2215
$ws.Range("F11").Value = 51
$ws.Range("F12").Value = 1676
$ws.Range("F13").Value = 2813
$ws.Range("F14").Value = 151
$ws.Range("F15").Value = 4190
$ws.Range("F16").Value = 355
$ws.Range("F17").Value = 183
$ws.Range("F19").Value = 534
$ws.Range("F20").Value = 249
$ws.Range("F21").Value = 40
$ws.Range("F23").Value = 93
$ws.Range("F24").Value = 287
$ws.Range("F25").Value = 4095
$ws.Range("F27").Value = 3577
$ws.Range("F28").Value = 1115
$ws.Range("F29").Value = 204
$ws.Range("F30").Value = 515
$ws.Range("F31").Value = 4357
$ws.Range("F32").Value = 78
$ws.Range("F33").Value = 391
$ws.Range("F34").Value = 459
$ws.Range("F35").Value = 368

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 3

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1016

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1016
$ws.Range("F8").Value = 5
$ws.Range("F10").Value = 42
$ws.Range("F11").Value = 725
$ws.Range("F12").Value = 2215
$ws.Range("F13").Value = 51
$ws.Range("F14").Value = 1676
$ws.Range("F16").Value = 2813
$ws.Range("F17").Value = 151
$ws.Range("F18").Value = 4190
$ws.Range("F19").Value = 355
$ws.Range("F20").Value = 183
$ws.Range("F22").Value = 534
$ws.Range("F23").Value = 249
$ws.Range("F24").Value = 40
$ws.Range("F27").Value = 93
$ws.Range("F28").Value = 287
$ws.Range("F29").Value = 4095
$ws.Range("F31").Value = 3577
$ws.Range("F32").Value = 1115
$ws.Range("F33").Value = 204
$ws.Range("F34").Value = 515
$ws.Range("F35").Value = 4357
$ws.Range("F36").Value = 78
$ws.Range("F37").Value = 391
$ws.Range("F38").Value = 459
$ws.Range("F39").Value = 368
$ws.Range("F40").Value = 3

Write-Output "Updated attendance counts (column F) across all 4 sheets"